$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-07-12 Friday" "2024-07-13 Saturday"

Replace-Text "560×2=1120" "822×9=7398"
Replace-Text "565×7=3955" "417×8=3336"
Replace-Text "727×6=4362" "761×3=2283"
Replace-Text "408×9=3672" "357×8=2856"
Replace-Text "897×8=7176" "666×9=5994"
Replace-Text "812×9=7308" "554×9=4986"
Replace-Text "573×7=4011" "902×5=4510"
Replace-Text "227×3=681" "599×7=4193"
Replace-Text "636×6=3816" "977×9=8793"
Replace-Text "769×6=4614" "910×6=5460"
Replace-Text "645×9=5805" "836×7=5852"
Replace-Text "655×6=3930" "831×3=2493"
Replace-Text "699×8=5592" "537×8=4296"
Replace-Text "976×9=8784" "889×2=1778"
Replace-Text "826×4=3304" "925×4=3700"
Replace-Text "811×5=4055" "965×2=1930"
Replace-Text "843×8=6744" "870×7=6090"
Replace-Text "163×8=1304" "310×2=620"
Replace-Text "522×3=1566" "310×6=1860"
Replace-Text "762×9=6858" "676×9=6084"
Replace-Text "583×5=2915" "521×4=2084"
Replace-Text "370×7=2590" "781×6=4686"
Replace-Text "172×2=344" "596×3=1788"
Replace-Text "797×3=2391" "236×2=472"
Replace-Text "170×7=1190" "403×8=3224"
